# Applies the "additional scraping" commit:
#  - adds a new "Player Info" sheet at the front
#  - adds a new "ODI Batting Extra" sheet at the end
#  - renames MATCH_CARD_LINK -> MATCH_CODE on the existing "ODI Batting" /
#    "ODI Bowling" sheets and replaces the full scorecard URL values with
#    just the numeric match code
#  - clears the stray empty INNING_NUMBER cell at ODI Batting!B3

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet, inserted before the first existing sheet
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
$col = 1
foreach ($h in $playerInfoHeaders) {
    $playerInfo.Cells.Item(1, $col).NumberFormat = "@"
    $playerInfo.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$playerInfoRow = @("4225", "Haris Sohail", "Left Handed", "Left Arm Orthodox")
$col = 1
foreach ($v in $playerInfoRow) {
    $playerInfo.Cells.Item(2, $col).NumberFormat = "@"
    $playerInfo.Cells.Item(2, $col).Value = $v
    $col = $col + 1
}

# ---------------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code,
#    drop the stray empty B3 cell
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")

$odiBatting.Range("D1").NumberFormat = "@"
$odiBatting.Range("D1").Value = "MATCH_CODE"

$battingMatchCodes = @("3533", "3535", "3538", "3546", "3713", "3715", "3717", "3719", "3720", "3745", "3747", "3751", "3757", "3770", "3772", "3789", "3792", "3797", "3798", "3799", "3801", "3802", "4110", "4114", "4200", "4227", "4273", "4274", "4275", "4276", "4277", "4287", "4292", "4294", "4304", "4334", "4337", "4340", "4349", "4375", "4376", "4432", "4686", "4688", "4690")

$r = 2
foreach ($code in $battingMatchCodes) {
    $odiBatting.Cells.Item($r, 4).NumberFormat = "@"
    $odiBatting.Cells.Item($r, 4).Value = $code
    $r = $r + 1
}

$odiBatting.Range("B3").Clear()

# ---------------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code
# ---------------------------------------------------------------------------
$odiBowling = $wb.Worksheets.Item("ODI Bowling")

$odiBowling.Range("B1").NumberFormat = "@"
$odiBowling.Range("B1").Value = "MATCH_CODE"

$bowlingMatchCodes = @("3713", "3715", "3717", "3719", "3720", "3745", "3747", "3751", "3757", "3772", "3789", "3792", "3797", "4110", "4200", "4274", "4275", "4276", "4277", "4292", "4294")

$r = 2
foreach ($code in $bowlingMatchCodes) {
    $odiBowling.Cells.Item($r, 2).NumberFormat = "@"
    $odiBowling.Cells.Item($r, 2).Value = $code
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4. New "ODI Batting Extra" sheet, appended after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
$col = 1
foreach ($h in $extraHeaders) {
    $battingExtra.Cells.Item(1, $col).NumberFormat = "@"
    $battingExtra.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$extraData = @(
    ,("4227", $null, $null, $null, $null, "NO")
    ,("4273", 3, "6", "1", "36.07%", "NO")
    ,("4274", 3, "5", "0", "11.97%", "NO")
    ,("4275", 3, "0", "0", "0.54%", "NO")
    ,("4276", $null, $null, $null, $null, "NO")
    ,("4277", 3, "11", "3", "42.35%", "NO")
    ,("4287", $null, $null, $null, $null, "NO")
    ,("4292", 5, "0", "0", "3.88%", "NO")
    ,("4294", $null, $null, $null, $null, "NO")
    ,("4304", 4, "1", "0", "7.62%", "NO")
    ,("4334", 5, "9", "3", "28.90%", "YES")
    ,("4337", $null, $null, $null, $null, "NO")
    ,("4340", 5, "2", "0", "11.74%", "NO")
    ,("4349", 5, "1", "0", "1.90%", "NO")
    ,("4375", 4, "1", "0", "13.11%", "NO")
    ,("4376", 5, "3", "1", "18.73%", "NO")
    ,("4432", $null, $null, $null, $null, "NO")
    ,("4686", $null, $null, $null, $null, "NO")
    ,("4688", 5, "0", "0", "5.49%", "NO")
    ,("4690", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraData) {
    $matchCode = $row[0]
    $battingPos = $row[1]
    $num4 = $row[2]
    $num6 = $row[3]
    $pct = $row[4]
    $mom = $row[5]

    $battingExtra.Cells.Item($r, 1).NumberFormat = "@"
    $battingExtra.Cells.Item($r, 1).Value = $matchCode

    if ($battingPos -ne $null) {
        $battingExtra.Cells.Item($r, 2).Value = $battingPos
    }

    if ($num4 -ne $null) {
        $battingExtra.Cells.Item($r, 3).NumberFormat = "@"
        $battingExtra.Cells.Item($r, 3).Value = $num4
    }

    if ($num6 -ne $null) {
        $battingExtra.Cells.Item($r, 4).NumberFormat = "@"
        $battingExtra.Cells.Item($r, 4).Value = $num6
    }

    if ($pct -ne $null) {
        $battingExtra.Cells.Item($r, 5).NumberFormat = "@"
        $battingExtra.Cells.Item($r, 5).Value = $pct
    }

    $battingExtra.Cells.Item($r, 6).NumberFormat = "@"
    $battingExtra.Cells.Item($r, 6).Value = $mom

    $r = $r + 1
}
